$d = $word.ActiveDocument

# 1. Remove " – Test as a Service" (en dash) from the title "TaaS – Test as a Service"
$d.Content.Find.Execute(" – Test as a Service", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Remove " - Test as a Service" (hyphen, bold) from the quoted text
$d.Content.Find.Execute(" - Test as a Service", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3. Remove the existing "_GoBack" bookmark (it wrapped the Architecture drawing).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 4. Re-add the "_GoBack" bookmark, collapsed, right after the "TaaS Server" bullet
#    item text (" Server"), before the paragraph mark. Directly adding a bookmark at
#    a collapsed Range sitting exactly on a paragraph-content boundary is unreliable,
#    so we insert a temporary unique marker run, anchor the bookmark just before it,
#    then delete the marker text (the now-empty bookmark stays behind in place).
$target = $d.Content
$target.Start = 0
$target.Find.Execute("TaaS Server", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0) | Out-Null
$target.InsertAfter("ZZZ_GOBACK_MARKER_ZZZ")

$markerRng = $d.Content
$markerRng.Start = 0
$markerRng.Find.Execute("ZZZ_GOBACK_MARKER_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRng.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRng)

$deleteRng = $d.Content
$deleteRng.Start = 0
$deleteRng.Find.Execute("ZZZ_GOBACK_MARKER_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteRng.Delete()

Write-Host "done"
